$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update T6 grade for Diego Neves Dos Santos (row 3) from blank to 1.25
$ws.Range("H3").Value = 1.25

# Update the active selection to H4 (matches the saved cursor position)
$ws.Range("H4").Select()
